{"js": "// Update the five 3-digit x 1-digit multiplication facts in each of the\n// first five populated rows of the practice table (25 cells total).\n// Each cell holds a single paragraph/run; we rewrite the run's text in\n// place (InsertLocation.replace) so the existing run/paragraph formatting\n// (font, size, justification) is preserved exactly as in the original.\n\nconst newValues = [\n  [\"560\u00d72=1120\", \"565\u00d77=3955\", \"727\u00d76=4362\", \"408\u00d79=3672\", \"897\u00d78=7176\"],\n  [\"812\u00d79=7308\", \"573\u00d77=4011\", \"227\u00d73=681\", \"636\u00d76=3816\", \"769\u00d76=4614\"],\n  [\"645\u00d79=5805\", \"655\u00d76=3930\", \"699\u00d78=5592\", \"976\u00d79=8784\", \"826\u00d74=3304\"],\n  [\"811\u00d75=4055\", \"843\u00d78=6744\", \"163\u00d78=1304\", \"522\u00d73=1566\", \"762\u00d79=6858\"],\n  [\"583\u00d75=2915\", \"370\u00d77=2590\", \"172\u00d72=344\", \"797\u00d73=2391\", \"170\u00d77=1190\"],\n];\n\n// The data rows, within the table, that actually contain text (the rest\n// are blank spacer rows).\nconst dataRowIndices = [0, 4, 9, 14, 19];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < dataRowIndices.length; r++) {\n  const rowIndex = dataRowIndices[r];\n  const values = newValues[r];\n  for (let c = 0; c < values.length; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n    const paragraph = paragraphs.items[0];\n    paragraph.insertText(values[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the five 3-digit x 1-digit multiplication facts in each of the\n# first five populated rows of the practice table (25 cells total).\n# Assigning Cell.Range.Text replaces just the text content of the cell's\n# existing run, preserving the run/paragraph formatting (font, size,\n# justification) exactly as in the original document.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based row indices of the table rows that actually hold data (the\n# remaining rows are blank spacer rows).\n$rowIndices = @(1, 5, 10, 15, 20)\n\n$values = @(\n    @(\"560\u00d72=1120\", \"565\u00d77=3955\", \"727\u00d76=4362\", \"408\u00d79=3672\", \"897\u00d78=7176\"),\n    @(\"812\u00d79=7308\", \"573\u00d77=4011\", \"227\u00d73=681\", \"636\u00d76=3816\", \"769\u00d76=4614\"),\n    @(\"645\u00d79=5805\", \"655\u00d76=3930\", \"699\u00d78=5592\", \"976\u00d79=8784\", \"826\u00d74=3304\"),\n    @(\"811\u00d75=4055\", \"843\u00d78=6744\", \"163\u00d78=1304\", \"522\u00d73=1566\", \"762\u00d79=6858\"),\n    @(\"583\u00d75=2915\", \"370\u00d77=2590\", \"172\u00d72=344\", \"797\u00d73=2391\", \"170\u00d77=1190\")\n)\n\nfor ($i = 0; $i -lt $rowIndices.Count; $i++) {\n    $rowIdx = $rowIndices[$i]\n    $rowVals = $values[$i]\n    for ($j = 0; $j -lt $rowVals.Count; $j++) {\n        $colIdx = $j + 1\n        $t.Cell($rowIdx, $colIdx).Range.Text = $rowVals[$j]\n    }\n}\n"}
